# Fix a typo in the first sheet's name ("moths" -> "months"), give that
# sheet previously-unset "natural" row heights for its first three rows,
# and make it the active/selected tab (it was previously on the
# "4 years" sheet).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Name = "6 months_៦ ខែ"

$ws.Rows.Item(1).RowHeight = 12.75
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 15

$ws.Activate()
